$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Structural change -----------------------------------------------------
# The original sheet had an unused blank row 1 above the header (header on
# row 2, data rows 3-24). The edited workbook removes that blank leading row
# (header moves to row 1, data starts on row 2) AND removes one data row -
# the "0% decay" (B=0) sample from the second ("Eu=50000000") group, which
# was a duplicate of the first group's B=0 row. Deleting both rows shifts
# everything else up accordingly and keeps all formulas/styles intact.

$ws.Rows(14).Delete()   # A=50000000, B=0 duplicate data row
$ws.Rows(1).Delete()    # blank leading row; header now on row 1

# --- Updated data values -----------------------------------------------------
# A handful of cells in the second ("Eu=50000000") group were re-measured /
# recalculated with the upload, independent of the row shift above. Values
# below are the new row numbers (after both deletions).

$ws.Cells.Item(13, 3).Value = 26.6
$ws.Cells.Item(13, 5).Value = 1.003

$ws.Cells.Item(14, 3).Value = 26.7
$ws.Cells.Item(14, 5).Value = 1.009
$ws.Cells.Item(14, 7).Value = 1

$ws.Cells.Item(15, 3).Value = 26.6
$ws.Cells.Item(15, 5).Value = 1.009
$ws.Cells.Item(15, 6).Value = 87.9

$ws.Cells.Item(16, 5).Value = 0.983
$ws.Cells.Item(16, 6).Value = 85.8
$ws.Cells.Item(16, 7).Value = 1.27

$ws.Cells.Item(17, 3).Value = 23.3
$ws.Cells.Item(17, 5).Value = 0.933
$ws.Cells.Item(17, 6).Value = 83.4

$ws.Cells.Item(18, 5).Value = 0.877
$ws.Cells.Item(18, 6).Value = 78.9
$ws.Cells.Item(18, 7).Value = 1.64

$ws.Cells.Item(19, 5).Value = 0.821
$ws.Cells.Item(19, 6).Value = 81.1

$ws.Cells.Item(20, 6).Value = 78.9
$ws.Cells.Item(20, 7).Value = 1.92

$ws.Cells.Item(21, 5).Value = 0.715
$ws.Cells.Item(21, 6).Value = 76.7

$ws.Cells.Item(22, 5).Value = 0.666
$ws.Cells.Item(22, 6).Value = 74.6

# --- Cosmetic view changes ---------------------------------------------------
$excel.ActiveWindow.Zoom = 175
$ws.Range("J14").Select()
